$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (008197156, Marcio, 46000) -> (004801481, Rafael, 51000)
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004801481"
$ws.Cells.Item(5, 1).ClearFormats()
$ws.Cells.Item(5, 2).Value = "Rafael"
$ws.Cells.Item(5, 3).Value = 51000

# Delete the old Rafael row (row 7: 004801481, Rafael, 35479.36)
$ws.Rows.Item(7).Delete()

# Delete the Ana row (originally row 9, now row 8 after the previous delete: 004267119, Ana, 13449.95)
$ws.Rows.Item(8).Delete()
